# Updates the "Recorded By" (column G) values on the 'Session Analysis Results'
# worksheet. For each listed row, the "System" token that used to lead the
# comma-separated list of recorder names is moved to the end of the list.
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#      "System, backup@backdoor.com, system" -> "backup@backdoor.com, system, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> expected old value / new value for column G ("Recorded By")
$rowUpdates = @{
    '2' = @{ Old = 'System, backup@backdoor.com, system'; New = 'backup@backdoor.com, system, System' };
    '3' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '5' = @{ Old = 'System, backup@backdoor.com'; New = 'backup@backdoor.com, System' };
    '6' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '8' = @{ Old = 'System, backup@backdoor.com'; New = 'backup@backdoor.com, System' };
    '10' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '11' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '12' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '13' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '14' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '15' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '17' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '18' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '19' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '20' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '21' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '22' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '24' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '26' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '28' = @{ Old = 'System, backup@backdoor.com, system'; New = 'backup@backdoor.com, system, System' };
    '29' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '31' = @{ Old = 'System, backup@backdoor.com'; New = 'backup@backdoor.com, System' };
    '32' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '34' = @{ Old = 'System, backup@backdoor.com'; New = 'backup@backdoor.com, System' };
    '36' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '37' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '38' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '39' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '40' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '41' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '43' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '44' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '45' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '46' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '47' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '48' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '50' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '52' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '54' = @{ Old = 'System, backup@backdoor.com, system'; New = 'backup@backdoor.com, system, System' };
    '55' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '57' = @{ Old = 'System, backup@backdoor.com'; New = 'backup@backdoor.com, System' };
    '58' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '60' = @{ Old = 'System, backup@backdoor.com'; New = 'backup@backdoor.com, System' };
    '62' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '63' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '64' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '65' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '66' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '67' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '69' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '70' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '71' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '72' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '73' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '74' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '76' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '78' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '80' = @{ Old = 'System, backup@backdoor.com'; New = 'backup@backdoor.com, System' };
    '81' = @{ Old = 'System, backup@backdoor.com'; New = 'backup@backdoor.com, System' };
    '82' = @{ Old = 'System, backup@backdoor.com'; New = 'backup@backdoor.com, System' };
    '83' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '84' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '85' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '86' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '90' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '92' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '93' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '94' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '96' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '99' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '101' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '106' = @{ Old = 'System, backup@backdoor.com'; New = 'backup@backdoor.com, System' };
    '107' = @{ Old = 'System, backup@backdoor.com'; New = 'backup@backdoor.com, System' };
    '108' = @{ Old = 'System, backup@backdoor.com'; New = 'backup@backdoor.com, System' };
    '109' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '110' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '111' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '112' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '116' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '118' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '119' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '120' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '122' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '125' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '127' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '132' = @{ Old = 'System, backup@backdoor.com'; New = 'backup@backdoor.com, System' };
    '133' = @{ Old = 'System, backup@backdoor.com'; New = 'backup@backdoor.com, System' };
    '134' = @{ Old = 'System, backup@backdoor.com'; New = 'backup@backdoor.com, System' };
    '135' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '136' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '137' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '138' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '142' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '144' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '145' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '146' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '148' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '151' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
    '153' = @{ Old = 'System, dnasr281@gmail.com'; New = 'dnasr281@gmail.com, System' };
}

foreach ($rowKey in $rowUpdates.Keys) {
    $row = [int]$rowKey
    $info = $rowUpdates[$rowKey]
    $cell = $ws.Cells.Item($row, 7)   # Column G = 7 ("Recorded By")

    $current = [string]$cell.Text
    if ($current -eq $info.Old) {
        $cell.Value = $info.New
    }
    else {
        # Fallback: value already updated or differs from expectation - apply
        # the transformation rule directly (move leading "System, " to the end)
        # so the edit is still applied even if text doesn't match exactly.
        if ($current -like 'System, *') {
            $rest = $current.Substring(8)
            $cell.Value = "$rest, System"
        }
    }
}

$wb.Save()
